$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the used range on column G and update cells whose value
# exactly matches one of the two known "Recorded By" strings, reordering
# the comma-separated author list as described in the commit diff.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row  # xlUp = -4162
if ($lastRow -lt 1) { $lastRow = 1 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # column G
    $val = $cell.Value2

    if ($val -eq "dnasr281@gmail.com, System") {
        $cell.Value = "System, dnasr281@gmail.com"
    }
    elseif ($val -eq "system, backup@backdoor.com, System") {
        $cell.Value = "backup@backdoor.com, system, System"
    }
}
